$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1444.6552
$ws.Range("I15").Value = 1444.6552
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4333.9656
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4164.9656

$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5350
$ws.Range("M40").ClearContents()

$ws.Range("H86").Value = 3102.3333
$ws.Range("I86").Value = 3161
$ws.Range("J86").Value = 3055.4
$ws.Range("K86").Value = 3161
$ws.Range("L86").Value = 3055.4
$ws.Range("M86").Value = -2038
$ws.Range("N86").Value = -5301.4

$ws.Range("H89").Value = 3102.3333
$ws.Range("I89").Value = 3161
$ws.Range("J89").Value = 3055.4
$ws.Range("K89").Value = 15805
$ws.Range("L89").Value = 15277
$ws.Range("M89").Value = -10189
$ws.Range("N89").Value = -26509

$ws.Range("H116").Value = 10833.333
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 10833.333
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 10833.333
$ws.Range("N116").Value = -17717.333
$ws.Range("M116").ClearContents()

$ws.Range("H138").Value = 2365.923
$ws.Range("I138").Value = 1852.4445
$ws.Range("J138").Value = 3521.25
$ws.Range("K138").Value = 5557.333500000001
$ws.Range("L138").Value = 10563.75
$ws.Range("M138").Value = -417.3335000000006
$ws.Range("N138").Value = -20843.75

$ws.Range("H141").Value = 102955.336
$ws.Range("I141").Value = 152864.67
$ws.Range("J141").Value = 3136.6667
$ws.Range("K141").Value = 458594.01
$ws.Range("L141").Value = 9410.000100000001
$ws.Range("M141").Value = -453414.01

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 37178
$ws.Range("I32").Value = 21579.55
$ws.Range("J32").Value = 132718.5
$ws.Range("K32").Value = 21579.55
$ws.Range("L32").Value = 132718.5
$ws.Range("M32").Value = -21292.55
$ws.Range("N32").Value = -133292.5

$ws.Range("H61").Value = 3828.1428
$ws.Range("I61").Value = 3386.875
$ws.Range("J61").Value = 4416.5
$ws.Range("K61").Value = 3386.875
$ws.Range("L61").Value = 4416.5
$ws.Range("M61").Value = -3174.875
$ws.Range("N61").Value = -4840.5

$ws.Range("H101").Value = 35601.75
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 35601.75
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 35601.75
$ws.Range("N101").Value = -42091.75

$ws.Range("H122").Value = 14903.3125
$ws.Range("I122").Value = 17759.691
$ws.Range("J122").Value = 2525.6667
$ws.Range("K122").Value = 53279.073
$ws.Range("L122").Value = 7577.000100000001
$ws.Range("M122").Value = -50829.073

$ws.Range("H132").Value = 2801.9697
$ws.Range("I132").Value = 2205.0527
$ws.Range("J132").Value = 3612.0715
$ws.Range("K132").Value = 6615.158100000001
$ws.Range("L132").Value = 10836.2145
$ws.Range("M132").Value = -4085.158100000001

$ws.Range("H136").Value = 3828.1428
$ws.Range("I136").Value = 3386.875
$ws.Range("J136").Value = 4416.5
$ws.Range("K136").Value = 10160.625
$ws.Range("L136").Value = 13249.5
$ws.Range("M136").Value = -7610.625
$ws.Range("N136").Value = -18349.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 124995
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 124995
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 124995
$ws.Range("N47").Value = -126035

$ws.Range("H82").Value = 19288.46
$ws.Range("I82").Value = 10154
$ws.Range("J82").Value = 24997.5
$ws.Range("K82").Value = 10154
$ws.Range("L82").Value = 24997.5
$ws.Range("M82").Value = -9771
$ws.Range("N82").Value = -25763.5

$ws.Range("H85").Value = 19288.46
$ws.Range("I85").Value = 10154
$ws.Range("J85").Value = 24997.5
$ws.Range("K85").Value = 10154
$ws.Range("L85").Value = 24997.5
$ws.Range("M85").Value = -8828
$ws.Range("N85").Value = -27649.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 10053
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 10053
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 10053
$ws.Range("N36").Value = -10829

$ws.Range("H40").Value = 10053
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10053
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 10053
$ws.Range("N40").Value = -10373

$ws.Range("H86").Value = 6503.5
$ws.Range("I86").Value = 3007
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 3007
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -1884
$ws.Range("N86").Value = -12246

$ws.Range("H89").Value = 6503.5
$ws.Range("I89").Value = 3007
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 15035
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -9419
$ws.Range("N89").Value = -61232

$ws.Range("H132").Value = 2149.0344
$ws.Range("I132").Value = 2175.28
$ws.Range("J132").Value = 1985
$ws.Range("K132").Value = 6525.84
$ws.Range("L132").Value = 5955
$ws.Range("M132").Value = -3995.84
$ws.Range("N132").Value = -11015

$ws.Range("H134").Value = 1906.4
$ws.Range("I134").Value = 1772.5518
$ws.Range("J134").Value = 2553.3333
$ws.Range("K134").Value = 5317.6554
$ws.Range("L134").Value = 7659.999899999999
$ws.Range("M134").Value = -2782.6554
$ws.Range("N134").Value = -12729.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3908.4167
$ws.Range("I118").Value = 230.5
$ws.Range("J118").Value = 5747.375
$ws.Range("K118").Value = 691.5
$ws.Range("L118").Value = 17242.125
$ws.Range("M118").Value = 551.5
$ws.Range("N118").Value = -19728.125

$ws.Range("H138").Value = 9247.25
$ws.Range("I138").Value = 3000
$ws.Range("J138").Value = 11329.667
$ws.Range("K138").Value = 9000
$ws.Range("L138").Value = 33989.001
$ws.Range("M138").Value = -3860
$ws.Range("N138").Value = -44269.001

$ws.Range("H139").Value = 5346.9165
$ws.Range("I139").Value = 3670.5
$ws.Range("J139").Value = 8699.75
$ws.Range("K139").Value = 11011.5
$ws.Range("L139").Value = 26099.25
$ws.Range("M139").Value = -5871.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3408000.8

$ws.Range("H80").Value = 8664.579
$ws.Range("I80").Value = 15506
$ws.Range("J80").Value = 3689
$ws.Range("K80").Value = 15506
$ws.Range("L80").Value = 3689
$ws.Range("M80").Value = -14508
$ws.Range("N80").Value = -5685

$ws.Range("H83").Value = 8664.579
$ws.Range("I83").Value = 15506
$ws.Range("J83").Value = 3689
$ws.Range("K83").Value = 77530
$ws.Range("L83").Value = 18445
$ws.Range("M83").Value = -72538
$ws.Range("N83").Value = -28429

$ws.Range("H113").Value = 83337336
$ws.Range("I113").Value = 166669330
$ws.Range("J113").Value = 5336
$ws.Range("K113").Value = 166669330
$ws.Range("L113").Value = 5336
$ws.Range("M113").Value = -166667160

$ws.Range("H122").Value = 1510.7858
$ws.Range("I122").Value = 1377.4546
$ws.Range("J122").Value = 1999.6666
$ws.Range("K122").Value = 4132.3638
$ws.Range("L122").Value = 5998.9998
$ws.Range("M122").Value = -1682.3638

$ws.Range("H132").Value = 3400
$ws.Range("I132").Value = 3400
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10200
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3056.1875
$ws.Range("I46").Value = 2379.8
$ws.Range("J46").Value = 3363.6365
$ws.Range("K46").Value = 2379.8
$ws.Range("L46").Value = 3363.6365
$ws.Range("M46").Value = -2191.8
$ws.Range("N46").Value = -3739.6365

$ws.Range("H132").Value = 2648.5818
$ws.Range("I132").Value = 2412.4048
$ws.Range("J132").Value = 3411.6155
$ws.Range("K132").Value = 7237.214399999999
$ws.Range("L132").Value = 10234.8465
$ws.Range("M132").Value = -4707.214399999999

$ws.Range("H136").Value = 3197.6956
$ws.Range("I136").Value = 2390.3
$ws.Range("J136").Value = 3818.7693
$ws.Range("K136").Value = 7170.900000000001
$ws.Range("L136").Value = 11456.3079
$ws.Range("M136").Value = -4620.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2023.92
$ws.Range("I122").Value = 1100.5
$ws.Range("J122").Value = 3665.5557
$ws.Range("K122").Value = 3301.5
$ws.Range("L122").Value = 10996.6671
$ws.Range("M122").Value = -851.5
$ws.Range("N122").Value = -15896.6671

$ws.Range("H132").Value = 8424.817999999999
$ws.Range("I132").Value = 9883.5625
$ws.Range("J132").Value = 4534.8335
$ws.Range("K132").Value = 29650.6875
$ws.Range("L132").Value = 13604.5005
$ws.Range("M132").Value = -27120.6875
